$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet to the generic default name
$ws.Name = "Sheet1"

# Clear the old "Mark" column values that no longer exist in the new dataset
$ws.Range("D2:D5").ClearContents()

# Replace the roster with 15 generic students (QE180001..QE180015 / Student 1..Student 15)
for ($i = 1; $i -le 15; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "A18C.DS"
}
for ($i = 1; $i -le 15; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = "QE18" + $i.ToString().PadLeft(4, '0')
}
for ($i = 1; $i -le 15; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = "Student " + $i
}
